$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J mirrors column I's formatting for rows 3-12 (new "2020" data column).
$ws.Range("I3:I12").Copy() | Out-Null
$ws.Range("J3:J12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New values for the 2020 column.
$ws.Range("J4").Value = 2020
$ws.Range("J5").Value = 253.27664777870578
$ws.Range("J7").Value = 93.236077839070575
$ws.Range("J8").Value = 160
$ws.Range("J10").Value = 69
$ws.Range("J11").Value = 48.5
$ws.Range("J12").Value = 22.8

# Move the active selection to the newly added header cell.
$ws.Range("J3").Select() | Out-Null
